$p = $ppt.ActivePresentation

# Slide 20: "Click on LimnoTech/GitHub-Training-SEMIFLD" -> add hyperlink
# to https://github.com/LimnoTech/GitHub-Training-SEMIFLD over the run(s)
# spanning "LimnoTech/GitHub-Training-SEMIFLD".
$s20 = $p.Slides.Item(20)
$tf20 = $s20.Shapes.Item(2).TextFrame
$tr20 = $tf20.TextRange
$para20 = $tr20.Paragraphs(5)
$linkRange = $para20.Characters(10, $para20.Length - 9)
$linkRange.ActionSettings.Item(1).Action = 7
$linkRange.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/LimnoTech/GitHub-Training-SEMIFLD"

# Slide 22: "Branch from Master" -> "Branch from Main"
$s22 = $p.Slides.Item(22)
$tf22 = $s22.Shapes.Item(2).TextFrame
$tr22 = $tf22.TextRange
$para22 = $tr22.Paragraphs(3)
$run22 = $para22.Runs(1)
$run22.Text = "Branch from Main"
